$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 13891371
$ws.Range("I100").Value = 34190236
$ws.Range("J100").Value = 2674.6316
$ws.Range("K100").Value = 34190236
$ws.Range("L100").Value = 2674.6316
$ws.Range("M100").Value = -34189695
$ws.Range("N100").Value = -3756.6316

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1451.7878
$ws.Range("I102").Value = 1450.3125
$ws.Range("J102").Value = 1499
$ws.Range("K102").Value = 1450.3125
$ws.Range("L102").Value = 1499
$ws.Range("M102").Value = 171.6875
$ws.Range("N102").Value = -4743

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20840680
$ws.Range("I20").Value = 26323812
$ws.Range("J20").Value = 4781.8
$ws.Range("K20").Value = 26323812
$ws.Range("L20").Value = 4781.8
$ws.Range("M20").Value = -26323565
$ws.Range("N20").Value = -5275.8
$ws.Range("H94").Value = 556.1875
$ws.Range("I94").Value = 538.38464
$ws.Range("J94").Value = 633.3333
$ws.Range("K94").Value = 538.38464
$ws.Range("L94").Value = 633.3333
$ws.Range("M94").Value = -87.38463999999999
$ws.Range("N94").Value = -1535.3333

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 55555876
$ws.Range("I22").Value = 100000080
$ws.Range("J22").Value = 623
$ws.Range("K22").Value = 100000080
$ws.Range("L22").Value = 623
$ws.Range("M22").Value = -99999730
$ws.Range("N22").Value = -1323
$ws.Range("H31").Value = 3396.1428
$ws.Range("I31").Value = 3345.4285
$ws.Range("J31").Value = 3497.5715
$ws.Range("K31").Value = 3345.4285
$ws.Range("L31").Value = 3497.5715
$ws.Range("M31").Value = -3050.4285
$ws.Range("N31").Value = -4087.5715
$ws.Range("H34").Value = 3396.1428
$ws.Range("I34").Value = 3345.4285
$ws.Range("J34").Value = 3497.5715
$ws.Range("K34").Value = 3345.4285
$ws.Range("L34").Value = 3497.5715
$ws.Range("M34").Value = -3143.4285
$ws.Range("N34").Value = -3901.5715
$ws.Range("H41").Value = 21000
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 31666.666
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 31666.666
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -32522.666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 67059.664
$ws.Range("I55").Value = 111759.43
$ws.Range("J55").Value = 4480
$ws.Range("K55").Value = 335278.29
$ws.Range("L55").Value = 13440
$ws.Range("M55").Value = -335101.29
$ws.Range("N55").Value = -13794
$ws.Range("H92").Value = 1106.3334
$ws.Range("J92").Value = 1387.875
$ws.Range("L92").Value = 4163.625
$ws.Range("N92").Value = -6659.625
$ws.Range("H132").Value = 1774.0526
$ws.Range("I132").Value = 1102.3334
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 9921.000599999999
$ws.Range("L132").Value = 17100
$ws.Range("M132").Value = -7391.000599999999
$ws.Range("N132").Value = -22160
$ws.Range("H139").Value = 2437.9678
$ws.Range("I139").Value = 983.8461
$ws.Range("K139").Value = 2951.5383
$ws.Range("M139").Value = 2188.4617

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 29000
$ws.Range("J32").Value = 29000
$ws.Range("L32").Value = 29000
$ws.Range("N32").Value = -29592
$ws.Range("H45").Value = 24000
$ws.Range("J45").Value = 24000
$ws.Range("L45").Value = 24000
$ws.Range("M45").Value = -25118
$ws.Range("H80").Value = 2325
$ws.Range("J80").Value = 2457.1428
$ws.Range("L80").Value = 2457.1428
$ws.Range("N80").Value = -4453.1428
$ws.Range("H83").Value = 2325
$ws.Range("J83").Value = 2457.1428
$ws.Range("L83").Value = 12285.714
$ws.Range("N83").Value = -22269.714
$ws.Range("H126").Value = 1240.8334
$ws.Range("I126").Value = 811.25
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 2433.75
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = 36.25
$ws.Range("N126").Value = -11240

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 529.5833
$ws.Range("I22").Value = 538.8823
$ws.Range("J22").Value = 507
$ws.Range("K22").Value = 538.8823
$ws.Range("L22").Value = 507
$ws.Range("M22").Value = -243.8823
$ws.Range("N22").Value = -1097
$ws.Range("H27").Value = 529.5833
$ws.Range("I27").Value = 538.8823
$ws.Range("J27").Value = 507
$ws.Range("K27").Value = 538.8823
$ws.Range("L27").Value = 507
$ws.Range("M27").Value = -431.8823
$ws.Range("N27").Value = -721
$ws.Range("H40").Value = 2019.65
$ws.Range("I40").Value = 2007.5264
$ws.Range("J40").Value = 2250
$ws.Range("K40").Value = 2007.5264
$ws.Range("L40").Value = 2250
$ws.Range("M40").Value = -1871.5264
$ws.Range("N40").Value = -2522

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9717
$ws.Range("I81").Value = 17554.5
$ws.Range("J81").Value = 2999.1428
$ws.Range("K81").Value = 35109
$ws.Range("L81").Value = 5998.2856
$ws.Range("M81").Value = -34048
$ws.Range("N81").Value = -8120.2856
$ws.Range("H84").Value = 9717
$ws.Range("I84").Value = 17554.5
$ws.Range("J84").Value = 2999.1428
$ws.Range("K84").Value = 175545
$ws.Range("L84").Value = 29991.428
$ws.Range("M84").Value = -170241
$ws.Range("N84").Value = -40599.428
$ws.Range("H87").Value = 39600
$ws.Range("I87").Value = 38000
$ws.Range("K87").Value = 38000
$ws.Range("M87").Value = -36752
$ws.Range("H90").Value = 39600
$ws.Range("I90").Value = 38000
$ws.Range("K90").Value = 114000
$ws.Range("M90").Value = -107760
$ws.Range("H96").Value = 2630.3076
$ws.Range("I96").Value = 2380
$ws.Range("J96").Value = 2786.75
$ws.Range("K96").Value = 2380
$ws.Range("L96").Value = 2786.75
$ws.Range("M96").Value = -1007
$ws.Range("N96").Value = -5532.75
$ws.Range("H100").Value = 3636827.5
$ws.Range("I100").Value = 5682274
$ws.Range("J100").Value = 477.77777
$ws.Range("K100").Value = 11364548
$ws.Range("L100").Value = 955.55554
$ws.Range("M100").Value = -11364007
$ws.Range("N100").Value = -2037.55554
$ws.Range("H107").Value = 5682.316
$ws.Range("I107").Value = 510.2857
$ws.Range("J107").Value = 8699.333000000001
$ws.Range("K107").Value = 1530.8571
$ws.Range("L107").Value = 26097.999
$ws.Range("M107").Value = 389.1428999999998
$ws.Range("N107").Value = -29937.999
$ws.Range("H122").Value = 4869.9
$ws.Range("I122").Value = 4528.4287
$ws.Range("K122").Value = 13585.2861
$ws.Range("M122").Value = -11135.2861
